$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 9423.714  # H51: 9425.429 -> 9423.714
$ws.Cells.Item(51, 9).Value = 8793.4  # I51: 8795.799999999999 -> 8793.4
$ws.Cells.Item(51, 11).Value = 8793.4  # K51: 8795.799999999999 -> 8793.4
$ws.Cells.Item(51, 13).Value = -8309.4  # M51: -8311.799999999999 -> -8309.4
$ws.Cells.Item(92, 8).Value = 50261.1  # H92: 50261.15 -> 50261.1
$ws.Cells.Item(92, 9).Value = 52801.156  # I92: 52801.21 -> 52801.156
$ws.Cells.Item(92, 11).Value = 52801.156  # K92: 52801.21 -> 52801.156
$ws.Cells.Item(92, 13).Value = -51553.156  # M92: -51553.21 -> -51553.156
$ws.Cells.Item(100, 8).Value = 1717.7778  # H100: 1719.4445 -> 1717.7778
$ws.Cells.Item(100, 9).Value = 1248.1333  # I100: 1250.1333 -> 1248.1333
$ws.Cells.Item(100, 11).Value = 1248.1333  # K100: 1250.1333 -> 1248.1333
$ws.Cells.Item(100, 13).Value = -707.1333  # M100: -709.1333 -> -707.1333
$ws.Cells.Item(107, 8).Value = 1262.6666  # H107: 1221.1428 -> 1262.6666
$ws.Cells.Item(107, 9).Value = 999.1111  # I107: 951.7895 -> 999.1111
$ws.Cells.Item(107, 11).Value = 999.1111  # K107: 951.7895 -> 999.1111
$ws.Cells.Item(107, 13).Value = 920.8889  # M107: 968.2105 -> 920.8889
$ws.Cells.Item(132, 8).Value = 3306.2766  # H132: 3462.7334 -> 3306.2766
$ws.Cells.Item(132, 9).Value = 2675.689  # I132: 2758 -> 2675.689
$ws.Cells.Item(132, 10).Value = 17494.5  # J132: 13329 -> 17494.5
$ws.Cells.Item(132, 11).Value = 8027.066999999999  # K132: 8274 -> 8027.066999999999
$ws.Cells.Item(132, 12).Value = 52483.5  # L132: 39987 -> 52483.5
$ws.Cells.Item(132, 13).Value = -5497.066999999999  # M132: -5744 -> -5497.066999999999
$ws.Cells.Item(132, 14).Value = -57543.5  # N132: -45047 -> -57543.5
$ws.Cells.Item(138, 8).Value = 5285.037  # H138: 5367.321 -> 5285.037
$ws.Cells.Item(138, 9).Value = 2923.5217  # I138: 3014.4092 -> 2923.5217
$ws.Cells.Item(138, 11).Value = 8770.5651  # K138: 9043.2276 -> 8770.5651
$ws.Cells.Item(138, 13).Value = -3630.5651  # M138: -3903.2276 -> -3630.5651
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 811.7  # H2: 619.8409 -> 811.7
$ws.Cells.Item(2, 9).Value = 647.5  # I2: 476.08334 -> 647.5
$ws.Cells.Item(2, 10).Value = 1263.25  # J2: 1266.75 -> 1263.25
$ws.Cells.Item(2, 11).Value = 647.5  # K2: 476.08334 -> 647.5
$ws.Cells.Item(2, 12).Value = 1263.25  # L2: 1266.75 -> 1263.25
$ws.Cells.Item(2, 13).Value = -534.5  # M2: -363.08334 -> -534.5
$ws.Cells.Item(2, 14).Value = -1489.25  # N2: -1492.75 -> -1489.25
$ws.Cells.Item(61, 8).Value = 1920.421  # H61: 2152.3157 -> 1920.421
$ws.Cells.Item(61, 9).Value = 1804.8889  # I61: 1806 -> 1804.8889
$ws.Cells.Item(61, 10).Value = 4000  # J61: 3999.3333 -> 4000
$ws.Cells.Item(61, 11).Value = 1804.8889  # K61: 1806 -> 1804.8889
$ws.Cells.Item(61, 12).Value = 4000  # L61: 3999.3333 -> 4000
$ws.Cells.Item(61, 13).Value = -1592.8889  # M61: -1594 -> -1592.8889
$ws.Cells.Item(61, 14).Value = -4424  # N61: -4423.3333 -> -4424
$ws.Cells.Item(74, 8).Value = 63861.875  # H74: 56897.277 -> 63861.875
$ws.Cells.Item(74, 9).Value = 84576.25  # I74: 78137.53999999999 -> 84576.25
$ws.Cells.Item(74, 10).Value = 1718.75  # J74: 1672.6 -> 1718.75
$ws.Cells.Item(74, 11).Value = 84576.25  # K74: 78137.53999999999 -> 84576.25
$ws.Cells.Item(74, 12).Value = 1718.75  # L74: 1672.6 -> 1718.75
$ws.Cells.Item(74, 13).Value = -83702.25  # M74: -77263.53999999999 -> -83702.25
$ws.Cells.Item(74, 14).Value = -3466.75  # N74: -3420.6 -> -3466.75
$ws.Cells.Item(77, 8).Value = 63861.875  # H77: 56897.277 -> 63861.875
$ws.Cells.Item(77, 9).Value = 84576.25  # I77: 78137.53999999999 -> 84576.25
$ws.Cells.Item(77, 10).Value = 1718.75  # J77: 1672.6 -> 1718.75
$ws.Cells.Item(77, 11).Value = 422881.25  # K77: 390687.7 -> 422881.25
$ws.Cells.Item(77, 12).Value = 8593.75  # L77: 8363 -> 8593.75
$ws.Cells.Item(77, 13).Value = -418513.25  # M77: -386319.7 -> -418513.25
$ws.Cells.Item(77, 14).Value = -17329.75  # N77: -17099 -> -17329.75
$ws.Cells.Item(116, 8).Value = 811.7  # H116: 619.8409 -> 811.7
$ws.Cells.Item(116, 9).Value = 647.5  # I116: 476.08334 -> 647.5
$ws.Cells.Item(116, 10).Value = 1263.25  # J116: 1266.75 -> 1263.25
$ws.Cells.Item(116, 11).Value = 647.5  # K116: 476.08334 -> 647.5
$ws.Cells.Item(116, 12).Value = 1263.25  # L116: 1266.75 -> 1263.25
$ws.Cells.Item(116, 13).Value = 1646.5  # M116: 1817.91666 -> 1646.5
$ws.Cells.Item(116, 14).Value = -5851.25  # N116: -5854.75 -> -5851.25
$ws.Cells.Item(132, 8).Value = 40986.816  # H132: 42651.848 -> 40986.816
$ws.Cells.Item(132, 9).Value = 55214.367  # I132: 61445.707 -> 55214.367
$ws.Cells.Item(132, 10).Value = 7196.375  # J132: 7152.3335 -> 7196.375
$ws.Cells.Item(132, 11).Value = 165643.101  # K132: 184337.121 -> 165643.101
$ws.Cells.Item(132, 12).Value = 21589.125  # L132: 21457.0005 -> 21589.125
$ws.Cells.Item(132, 13).Value = -163113.101  # M132: -181807.121 -> -163113.101
$ws.Cells.Item(132, 14).Value = -26649.125  # N132: -26517.0005 -> -26649.125
$ws.Cells.Item(136, 8).Value = 1920.421  # H136: 2152.3157 -> 1920.421
$ws.Cells.Item(136, 9).Value = 1804.8889  # I136: 1806 -> 1804.8889
$ws.Cells.Item(136, 10).Value = 4000  # J136: 3999.3333 -> 4000
$ws.Cells.Item(136, 11).Value = 5414.6667  # K136: 5418 -> 5414.6667
$ws.Cells.Item(136, 12).Value = 12000  # L136: 11997.9999 -> 12000
$ws.Cells.Item(136, 13).Value = -2864.6667  # M136: -2868 -> -2864.6667
$ws.Cells.Item(136, 14).Value = -17100  # N136: -17097.9999 -> -17100
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 811.7  # H3: 619.8409 -> 811.7
$ws.Cells.Item(3, 9).Value = 647.5  # I3: 476.08334 -> 647.5
$ws.Cells.Item(3, 10).Value = 1263.25  # J3: 1266.75 -> 1263.25
$ws.Cells.Item(3, 11).Value = 647.5  # K3: 476.08334 -> 647.5
$ws.Cells.Item(3, 12).Value = 1263.25  # L3: 1266.75 -> 1263.25
$ws.Cells.Item(3, 13).Value = -533.5  # M3: -362.08334 -> -533.5
$ws.Cells.Item(3, 14).Value = -1491.25  # N3: -1494.75 -> -1491.25
$ws.Cells.Item(105, 8).Value = 4504.72  # H105: 4383.1113 -> 4504.72
$ws.Cells.Item(105, 9).Value = 4331.25  # I105: 4159.619 -> 4331.25
$ws.Cells.Item(105, 10).Value = 5198.6  # J105: 5165.3335 -> 5198.6
$ws.Cells.Item(105, 11).Value = 4331.25  # K105: 4159.619 -> 4331.25
$ws.Cells.Item(105, 12).Value = 5198.6  # L105: 5165.3335 -> 5198.6
$ws.Cells.Item(105, 13).Value = -2584.25  # M105: -2412.619 -> -2584.25
$ws.Cells.Item(105, 14).Value = -8692.6  # N105: -8659.333500000001 -> -8692.6
$ws.Cells.Item(134, 8).Value = 1719.0454  # H134: 1847.9048 -> 1719.0454
$ws.Cells.Item(134, 9).Value = 1719.0454  # I134: 1792.5 -> 1719.0454
$ws.Cells.Item(134, 10).Value = 0  # J134: 2956 -> 0
$ws.Cells.Item(134, 11).Value = 5157.1362  # K134: 5377.5 -> 5157.1362
$ws.Cells.Item(134, 12).Value = 0  # L134: 8868 -> 0
$ws.Cells.Item(134, 13).Value = -2622.1362  # M134: -2842.5 -> -2622.1362
$ws.Cells.Item(134, 14).Value = ""  # remove N134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 799.3333  # H16: 818.05884 -> 799.3333
$ws.Cells.Item(16, 9).Value = 758.1177  # I16: 775.4375 -> 758.1177
$ws.Cells.Item(16, 11).Value = 758.1177  # K16: 775.4375 -> 758.1177
$ws.Cells.Item(16, 13).Value = -471.1177  # M16: -488.4375 -> -471.1177
$ws.Cells.Item(105, 8).Value = 804  # H105: 753.7 -> 804
$ws.Cells.Item(105, 9).Value = 804  # I105: 753.7 -> 804
$ws.Cells.Item(105, 11).Value = 804  # K105: 753.7 -> 804
$ws.Cells.Item(105, 13).Value = 943  # M105: 993.3 -> 943
$ws.Cells.Item(113, 8).Value = 799.3333  # H113: 818.05884 -> 799.3333
$ws.Cells.Item(113, 9).Value = 758.1177  # I113: 775.4375 -> 758.1177
$ws.Cells.Item(113, 11).Value = 758.1177  # K113: 775.4375 -> 758.1177
$ws.Cells.Item(113, 13).Value = 1411.8823  # M113: 1394.5625 -> 1411.8823
$ws.Cells.Item(133, 8).Value = 99816  # H133: 99859.2 -> 99816
$ws.Cells.Item(133, 9).Value = 99699  # I133: 99799 -> 99699
$ws.Cells.Item(133, 10).Value = 99839.39999999999  # J133: 99874.25 -> 99839.39999999999
$ws.Cells.Item(133, 11).Value = 99699  # K133: 99799 -> 99699
$ws.Cells.Item(133, 12).Value = 99839.39999999999  # L133: 99874.25 -> 99839.39999999999
$ws.Cells.Item(133, 13).Value = -97169  # M133: -97269 -> -97169
$ws.Cells.Item(133, 14).Value = -104899.4  # N133: -104934.25 -> -104899.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 103.4  # H12: 21.8 -> 103.4
$ws.Cells.Item(12, 9).Value = 250.5  # I12: 26 -> 250.5
$ws.Cells.Item(12, 10).Value = 5.3333335  # J12: 15.5 -> 5.3333335
$ws.Cells.Item(12, 11).Value = 751.5  # K12: 78 -> 751.5
$ws.Cells.Item(12, 12).Value = 16.0000005  # L12: 46.5 -> 16.0000005
$ws.Cells.Item(12, 13).Value = -578.5  # M12: 95 -> -578.5
$ws.Cells.Item(12, 14).Value = -362.0000005  # N12: -392.5 -> -362.0000005
$ws.Cells.Item(44, 8).Value = 523.9259  # H44: 805.2 -> 523.9259
$ws.Cells.Item(44, 9).Value = 381.08334  # I44: 506.5 -> 381.08334
$ws.Cells.Item(44, 10).Value = 1666.6666  # J44: 2000 -> 1666.6666
$ws.Cells.Item(44, 11).Value = 1143.25002  # K44: 1519.5 -> 1143.25002
$ws.Cells.Item(44, 12).Value = 4999.9998  # L44: 6000 -> 4999.9998
$ws.Cells.Item(44, 13).Value = -745.2500199999999  # M44: -1121.5 -> -745.2500199999999
$ws.Cells.Item(44, 14).Value = -5795.9998  # N44: -6796 -> -5795.9998
$ws.Cells.Item(109, 8).Value = 1001182.8  # H109: 1169.6666 -> 1001182.8
$ws.Cells.Item(109, 9).Value = 1166  # I109: 1003.375 -> 1166
$ws.Cells.Item(109, 10).Value = 5001250  # J109: 2500 -> 5001250
$ws.Cells.Item(109, 11).Value = 3498  # K109: 3010.125 -> 3498
$ws.Cells.Item(109, 12).Value = 15003750  # L109: 7500 -> 15003750
$ws.Cells.Item(109, 13).Value = -2458  # M109: -1970.125 -> -2458
$ws.Cells.Item(109, 14).Value = -15005830  # N109: -9580 -> -15005830
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3172.2307  # H80: 3183.6924 -> 3172.2307
$ws.Cells.Item(80, 9).Value = 2931.2856  # I80: 2871 -> 2931.2856
$ws.Cells.Item(80, 10).Value = 3453.3333  # J80: 3684 -> 3453.3333
$ws.Cells.Item(80, 11).Value = 2931.2856  # K80: 2871 -> 2931.2856
$ws.Cells.Item(80, 12).Value = 3453.3333  # L80: 3684 -> 3453.3333
$ws.Cells.Item(80, 13).Value = -1933.2856  # M80: -1873 -> -1933.2856
$ws.Cells.Item(80, 14).Value = -5449.3333  # N80: -5680 -> -5449.3333
$ws.Cells.Item(83, 8).Value = 3172.2307  # H83: 3183.6924 -> 3172.2307
$ws.Cells.Item(83, 9).Value = 2931.2856  # I83: 2871 -> 2931.2856
$ws.Cells.Item(83, 10).Value = 3453.3333  # J83: 3684 -> 3453.3333
$ws.Cells.Item(83, 11).Value = 14656.428  # K83: 14355 -> 14656.428
$ws.Cells.Item(83, 12).Value = 17266.6665  # L83: 18420 -> 17266.6665
$ws.Cells.Item(83, 13).Value = -9664.428  # M83: -9363 -> -9664.428
$ws.Cells.Item(83, 14).Value = -27250.6665  # N83: -28404 -> -27250.6665
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 36892.18  # H22: 36894.18 -> 36892.18
$ws.Cells.Item(22, 9).Value = 45999.227  # I22: 46001.773 -> 45999.227
$ws.Cells.Item(22, 11).Value = 45999.227  # K22: 46001.773 -> 45999.227
$ws.Cells.Item(22, 13).Value = -45704.227  # M22: -45706.773 -> -45704.227
$ws.Cells.Item(27, 8).Value = 36892.18  # H27: 36894.18 -> 36892.18
$ws.Cells.Item(27, 9).Value = 45999.227  # I27: 46001.773 -> 45999.227
$ws.Cells.Item(27, 11).Value = 45999.227  # K27: 46001.773 -> 45999.227
$ws.Cells.Item(27, 13).Value = -45892.227  # M27: -45894.773 -> -45892.227
$ws.Cells.Item(82, 8).Value = 2251.4119  # H82: 2253.2058 -> 2251.4119
$ws.Cells.Item(82, 9).Value = 1318.2  # I82: 1327.3 -> 1318.2
$ws.Cells.Item(82, 10).Value = 2640.25  # J82: 2639 -> 2640.25
$ws.Cells.Item(82, 11).Value = 1318.2  # K82: 1327.3 -> 1318.2
$ws.Cells.Item(82, 12).Value = 2640.25  # L82: 2639 -> 2640.25
$ws.Cells.Item(82, 13).Value = -957.2  # M82: -966.3 -> -957.2
$ws.Cells.Item(82, 14).Value = -3362.25  # N82: -3361 -> -3362.25
$ws.Cells.Item(85, 8).Value = 2251.4119  # H85: 2253.2058 -> 2251.4119
$ws.Cells.Item(85, 9).Value = 1318.2  # I85: 1327.3 -> 1318.2
$ws.Cells.Item(85, 10).Value = 2640.25  # J85: 2639 -> 2640.25
$ws.Cells.Item(85, 11).Value = 1318.2  # K85: 1327.3 -> 1318.2
$ws.Cells.Item(85, 12).Value = 2640.25  # L85: 2639 -> 2640.25
$ws.Cells.Item(85, 13).Value = -70.20000000000005  # M85: -79.29999999999995 -> -70.20000000000005
$ws.Cells.Item(85, 14).Value = -5136.25  # N85: -5135 -> -5136.25
$ws.Cells.Item(132, 8).Value = 40061.47  # H132: 41321.547 -> 40061.47
$ws.Cells.Item(132, 9).Value = 45102.75  # I132: 46736.223 -> 45102.75
$ws.Cells.Item(132, 11).Value = 135308.25  # K132: 140208.669 -> 135308.25
$ws.Cells.Item(132, 13).Value = -132778.25  # M132: -137678.669 -> -132778.25
$ws.Cells.Item(136, 8).Value = 4642.7896  # H136: 4689.294 -> 4642.7896
$ws.Cells.Item(136, 9).Value = 3846.5833  # I136: 3866.4 -> 3846.5833
$ws.Cells.Item(136, 10).Value = 6007.7144  # J136: 5864.857 -> 6007.7144
$ws.Cells.Item(136, 11).Value = 11539.7499  # K136: 11599.2 -> 11539.7499
$ws.Cells.Item(136, 12).Value = 18023.1432  # L136: 17594.571 -> 18023.1432
$ws.Cells.Item(136, 13).Value = -8989.749899999999  # M136: -9049.200000000001 -> -8989.749899999999
$ws.Cells.Item(136, 14).Value = -23123.1432  # N136: -22694.571 -> -23123.1432
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 8866.5  # H4: 8783 -> 8866.5
$ws.Cells.Item(4, 9).Value = 23349.5  # I4: 15732.667 -> 23349.5
$ws.Cells.Item(4, 10).Value = 1625  # J4: 1833.3334 -> 1625
$ws.Cells.Item(4, 11).Value = 23349.5  # K4: 15732.667 -> 23349.5
$ws.Cells.Item(4, 12).Value = 1625  # L4: 1833.3334 -> 1625
$ws.Cells.Item(4, 13).Value = -23236.5  # M4: -15619.667 -> -23236.5
$ws.Cells.Item(4, 14).Value = -1851  # N4: -2059.3334 -> -1851
$ws.Cells.Item(99, 8).Value = 98476  # H99: 0 -> 98476
$ws.Cells.Item(99, 10).Value = 98476  # J99: 0 -> 98476
$ws.Cells.Item(99, 12).Value = 98476  # L99: 0 -> 98476
$ws.Cells.Item(99, 14).Value = -104466  # N99: ADD -> -104466
$ws.Cells.Item(133, 8).Value = 89999  # H133: 89994.5 -> 89999
$ws.Cells.Item(133, 10).Value = 89999  # J133: 89994.5 -> 89999
$ws.Cells.Item(133, 12).Value = 89999  # L133: 89994.5 -> 89999
$ws.Cells.Item(133, 14).Value = -100119  # N133: -100114.5 -> -100119
